# Updates the crypto price table to the "Tue Dec 13 20:10:30 UTC 2022" GitHub Actions
# snapshot (see commit message). Every row keeps its refreshed Price (column D) and
# Hora/hour stamp (column G, 19 -> 20); two coin pairs (rows 42/43 and 49/50) also swap
# places in the ranking, so their Coin/Link/Price/Volume cells trade contents too.
#
# All of these cells hold TEXT (e.g. "268.36", "20"), not numbers, even though most
# values look numeric. Writing a numeric-looking string straight into Range.Value lets
# Excel auto-convert it to a real number, so instead we stage the text in an unused
# helper cell (forcing text with a leading apostrophe) and use
# PasteSpecial(xlPasteValues) to copy only the literal text into the destination cell,
# leaving that cell's existing number format/style untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$helper = $ws.Range("J1")
$xlPasteValues = -4163

function Set-TextValue($cellRef, $text) {
    $helper.Value = "'" + $text
    $helper.Copy() | Out-Null
    $ws.Range($cellRef).PasteSpecial($xlPasteValues)
}


# Row 2
Set-TextValue "D2" "268.36"
Set-TextValue "G2" "20"

# Row 3
Set-TextValue "D3" "23.00"
Set-TextValue "G3" "20"

# Row 4
Set-TextValue "D4" "6.341"
Set-TextValue "G4" "20"

# Row 5
Set-TextValue "D5" "0.06166"
Set-TextValue "G5" "20"

# Row 6
Set-TextValue "D6" "3.658"
Set-TextValue "G6" "20"

# Row 7
Set-TextValue "D7" "6.688"
Set-TextValue "G7" "20"

# Row 8
Set-TextValue "D8" "1.387"
Set-TextValue "G8" "20"

# Row 9
Set-TextValue "D9" "0.8299"
Set-TextValue "G9" "20"

# Row 10
Set-TextValue "D10" "0.01372"
Set-TextValue "G10" "20"

# Row 11
Set-TextValue "D11" "0.1605"
Set-TextValue "G11" "20"

# Row 12
Set-TextValue "D12" "0.08277"
Set-TextValue "G12" "20"

# Row 13
Set-TextValue "D13" "0.03485"
Set-TextValue "G13" "20"

# Row 14
Set-TextValue "D14" "0.03199"
Set-TextValue "G14" "20"

# Row 15
Set-TextValue "D15" "0.09326"
Set-TextValue "G15" "20"

# Row 16
Set-TextValue "D16" "3.845"
Set-TextValue "G16" "20"

# Row 17
Set-TextValue "D17" "0.001638"
Set-TextValue "G17" "20"

# Row 18
Set-TextValue "D18" "0.04740"
Set-TextValue "G18" "20"

# Row 19
Set-TextValue "D19" "0.006349"
Set-TextValue "G19" "20"

# Row 20
Set-TextValue "D20" "0.005652"
Set-TextValue "G20" "20"

# Row 21
Set-TextValue "D21" "0.001078"
Set-TextValue "G21" "20"

# Row 22
Set-TextValue "D22" "0.0001501"
Set-TextValue "G22" "20"

# Row 23
Set-TextValue "G23" "20"

# Row 24
Set-TextValue "G24" "20"

# Row 25
Set-TextValue "G25" "20"

# Row 26
Set-TextValue "G26" "20"

# Row 27
Set-TextValue "D27" "0.0002705"
Set-TextValue "G27" "20"

# Row 28
Set-TextValue "G28" "20"

# Row 29
Set-TextValue "G29" "20"

# Row 30
Set-TextValue "G30" "20"

# Row 31
Set-TextValue "G31" "20"

# Row 32
Set-TextValue "G32" "20"

# Row 33
Set-TextValue "G33" "20"

# Row 34
Set-TextValue "G34" "20"

# Row 35
Set-TextValue "G35" "20"

# Row 36
Set-TextValue "G36" "20"

# Row 37
Set-TextValue "G37" "20"

# Row 38
Set-TextValue "G38" "20"

# Row 39
Set-TextValue "G39" "20"

# Row 40
Set-TextValue "D40" "0.04703"
Set-TextValue "G40" "20"

# Row 41
Set-TextValue "D41" "0.006960"
Set-TextValue "G41" "20"

# Row 42
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D42" "0.003802"
$ws.Range("E42").Value = "41CEJICEJI"
Set-TextValue "G42" "20"

# Row 43
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D43" "0.1157"
$ws.Range("E43").Value = "42BKEXTokenBKK"
Set-TextValue "G43" "20"

# Row 44
Set-TextValue "D44" "0.01154"
Set-TextValue "G44" "20"

# Row 45
Set-TextValue "D45" "0.00006254"
Set-TextValue "G45" "20"

# Row 46
Set-TextValue "D46" "0.0009905"
Set-TextValue "G46" "20"

# Row 47
Set-TextValue "G47" "20"

# Row 48
Set-TextValue "D48" "0.9205"
Set-TextValue "G48" "20"

# Row 49
$ws.Range("B49").Value = "BOLO"
$ws.Range("C49").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
Set-TextValue "D49" "0.002249"
$ws.Range("E49").Value = "48BOLOBOLO"
Set-TextValue "G49" "20"

# Row 50
$ws.Range("B50").Value = "CryptobidCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/h39bvStAP+cryptobidcoin-cbc"
Set-TextValue "D50" "0.00001401"
$ws.Range("E50").Value = "49CryptobidCoinCBCWorstin24h"
Set-TextValue "G50" "20"

# Row 51
Set-TextValue "D51" "0.01241"
Set-TextValue "G51" "20"

$helper.Clear() | Out-Null
$excel.CutCopyMode = $false

